$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Replace "Harmonized-DB" with "Harmonized" in the Variable column (D2:D4)
$ws.Range("D2").Value = "prefix|Emissions|BC|Harmonized"
$ws.Range("D3").Value = "prefix|Emissions|BC|sector1|Harmonized"
$ws.Range("D4").Value = "prefix|Emissions|BC|sector2|Harmonized"

# Update the selection to match the saved view state
$ws.Range("L9").Select()
